$d = $word.ActiveDocument

# 1. Remove the standalone "в число студентов " run (with its white
#    highlight/shading) from the main body text. The run sits between
#    ", " and the "{" that opens the {levelEducation} merge field, so
#    deleting just this run's text collapses ", " + "{levelEducation}"
#    back together without touching anything else.
$bodyRange = $d.Content
$foundBody = $bodyRange.Find.Execute("в число студентов ", $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)
if ($foundBody) {
    $bodyRange.Delete()
}

# 2. Collapse the "от «___»___2022 №___" header line back into a single
#    run (it was previously split across three runs with grammar-check
#    markers — <w:proofErr w:type="gramStart"/>/<w:proofErr w:type="gramEnd"/>
#    — sitting between them). Re-typing the same visible text over itself
#    via Find & Replace merges it into one run and drops the stale
#    proofErr markers. Search every header of every section since the
#    text lives in whichever header holds the order-approval line.
$sections = $d.Sections
for ($s = 1; $s -le $sections.Count; $s++) {
    $section = $sections.Item($s)
    $headers = $section.Headers
    for ($i = 1; $i -le $headers.Count; $i++) {
        $h = $headers.Item($i)
        if ($h.Exists) {
            $hRange = $h.Range
            [void]$hRange.Find.Execute("от «___»___2022 №___", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "от «___»___2022 №___", 2)
        }
    }
}
